# Day 5 Sprint backlog 29/01/21
# Applies the "Day 5" updates to the Sprint 1 backlog worksheet:
#  - Gina also now on the "send questionnaires" story (F16)
#  - Row 20 "Set up database stuff" task marked Done on 29/01/21 (F20, new cell)
#  - Wording tweak to the daily-scrum note (I34)
#  - New Scrum-minutes row for day 5 / 29-01-21 (row 30)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Small text edits -------------------------------------------------

# "Gina" -> "Gina, Shaun"
$ws.Range("F16").Value = "Gina, Shaun"

# Daily scrum note: drop "to solve problems"
$ws.Range("I34").Value = "daily scrum shouldn't cover the technicalities of work, it should focus on the next steps and help keep things moving - technicalities can be done at other times aside from the scrum "

# --- New "Done?" date for the database task (row 20) -------------------

# F20 didn't exist before; give it the same look as the other "Done?" date
# cells in that mini-table (e.g. F21 = "Yes! 27/01/21") by copying F21's
# formatting, then filling in the new value.
$ws.Range("F21").Copy()
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("F20").Value = "Yes! 29/01/21"

# --- New Scrum Minutes row for day 5 (29/01/21), row 30 ----------------

# Date cell: copy the date-formatted style from the row above (A29) so the
# new date renders the same way (dd/mm/yy).
$ws.Range("A29").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A30").Value = 44225

# Text cells: copy formatting from the equivalent cells in row 29.
$ws.Range("B29").Copy()
$ws.Range("B30").PasteSpecial(-4122)
$ws.Range("B30").Value = "1 B, 5 S"

$ws.Range("B29").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = "- Yesterday hard, SQL linked now"

$ws.Range("E29").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = "- Feeling Lost                                                  - Left behind                                - File architecture confusing"

$ws.Range("G29").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G30").Value = "- Get all made pages working                                               - Get some basic SQL function to show off                                                                - Sprint Review"

# Merge the split label cells the same way the other scrum-minute rows do.
$ws.Range("C30:D30").Merge()
$ws.Range("E30:F30").Merge()
